$d = $word.ActiveDocument

# The "Date:" field in the signature block reads "09/10/2021"; change the
# day from "10" to "12" so it reads "09/12/2021".
$d.Content.Find.Execute("09/10/2021", $true, $false, $false, $false, $false,
                         $true, 1, $false, "09/12/2021", 2)
